# "Foi ajustado os dias de cada sprint." — fill in the start/end dates that
# were still blank for each sprint's date range, and leave the selection
# where the user last clicked (the Sprint 2 header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 1 (rows 3-7, merged G3:G7 / H3:H7): start date was already set
# (26/04/2015); the end date had been left empty - fill it in now.
$ws.Range("H3").Value = 42131          # 07/05/2015
$ws.Range("H3").NumberFormat = $ws.Range("G3").NumberFormat

# Sprint 2 (rows 10-15, merged G10:G15 / H10:H15): neither date had been
# filled in yet.
$ws.Range("G10").Value = 42132         # 08/05/2015
$ws.Range("G10").NumberFormat = $ws.Range("G3").NumberFormat
$ws.Range("H10").Value = 42143         # 19/05/2015
$ws.Range("H10").NumberFormat = $ws.Range("G3").NumberFormat

# Sprint 3 (rows 18-23, merged G18:G23 / H18:H23): neither date had been
# filled in yet.
$ws.Range("G18").Value = 42144         # 20/05/2015
$ws.Range("G18").NumberFormat = $ws.Range("G3").NumberFormat
$ws.Range("H18").Value = 42155         # 31/05/2015
$ws.Range("H18").NumberFormat = $ws.Range("G3").NumberFormat

# Leave the selection on the Sprint 2 title bar, matching where the user
# ended up after editing.
$ws.Range("A8:H8").Select()
